$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "FT23250108DDM2M9"
